$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "3) GUI Design weiter"
$ws.Range("C5").Value = 3
$ws.Range("E10").Select()
